# "Add thêm nhân sự Nguyễn Hữu Quang" — a new personal order (HD-LUXURY 707)
# is added to the "Đơn 1 bác sĩ" sheet, which shifts the existing "Tổng"
# (total) row down one row and updates its totals; the dependent payroll
# figures on the "Lương" sheet are updated to match.

$wb = $excel.ActiveWorkbook

# --- Sheet "Đơn 1 bác sĩ": insert a new order row before the "Tổng" row ---
$ws = $wb.Worksheets.Item("Đơn 1 bác sĩ")

# Push the existing "Tổng" row (row 3) down to row 4 to make room for the
# new order, then fill in the new row 3 with the new order's data.
$ws.Rows.Item(3).Insert()

$ws.Cells.Item(3, 1).Value = "HD-LUXURY"
$ws.Cells.Item(3, 2).Value = 707

# Keep the execution date as literal text (matches the existing "Ngày thực
# hiện" column, which stores dates as plain text, not Excel date serials).
$ws.Cells.Item(3, 3).NumberFormat = "@"
$ws.Cells.Item(3, 3).Value = "08-31-2024"
$ws.Cells.Item(3, 3).ClearFormats()

$ws.Cells.Item(3, 4).Value = "SÓC TRĂNG"
$ws.Cells.Item(3, 5).Value = "dương ngọc hân"
$ws.Cells.Item(3, 6).Value = "Cá nhân"
$ws.Cells.Item(3, 7).Value = "Tiêm Filler"
$ws.Cells.Item(3, 8).Value = 1300000
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 1300000
$ws.Cells.Item(3, 12).Value = 1300000
$ws.Cells.Item(3, 13).Value = 0.1
$ws.Cells.Item(3, 14).Value = 130000

# Update the "Tổng" (Total) row, now on row 4, to include the new order.
$ws.Cells.Item(4, 2).Value = 2
$ws.Cells.Item(4, 8).Value = 5300000
$ws.Cells.Item(4, 11).Value = 5300000
$ws.Cells.Item(4, 12).Value = 5300000
$ws.Cells.Item(4, 14).Value = 530000

# --- Sheet "Lương": update payroll figures affected by the new order ---
$wsL = $wb.Worksheets.Item("Lương")

$wsL.Cells.Item(1, 2).Value = 16
$wsL.Cells.Item(22, 2).Value = 24.5
$wsL.Cells.Item(23, 2).Value = 857500
$wsL.Cells.Item(24, 2).Value = 4812500
$wsL.Cells.Item(27, 2).Value = 530000
$wsL.Cells.Item(34, 2).Value = 4800000
$wsL.Cells.Item(35, 2).Value = 4800000
